# [IMP] midea backport from 12.0
# Adds a new "G/C RC" (Goods Clearing Reconciliation) account row to the
# account.account_test chart-of-accounts sheet, just above the existing
# "Undistributed Profits/Losses" placeholder row, and updates the sheet's
# view state (scroll/selection) to point at the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a brand-new row right before the last row (old row 38, the
#    "z0bug.lp" / Undistributed Profits/Losses row). Excel shifts that
#    row (and everything below it, nothing here) down to row 39 while
#    inheriting the formatting of the row above for the newly inserted
#    row 38.
# ---------------------------------------------------------------------
$ws.Rows.Item(38).Insert()

# ---------------------------------------------------------------------
# 2. Populate the new row 38 with the new G/C RC account data:
#      id            -> z0bug.coa_gc_rc
#      code          -> 490050
#      name          -> G/C RC
#      user_type_id  -> account.data_account_type_current_liabilities
#      reconcile     -> False
# ---------------------------------------------------------------------
$ws.Range("A38").Value = "z0bug.coa_gc_rc"
$ws.Range("B38").Value = 490050
$ws.Range("C38").Value = "G/C RC"
$ws.Range("D38").Value = "account.data_account_type_current_liabilities"

# Copy the "False" text value from an existing reconcile cell (E2) so the
# new E38 cell matches the same text type/format instead of being
# auto-detected as a boolean.
$ws.Range("E2").Copy()
$ws.Range("E38").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Update the view state: keep row 1 frozen as the header, scroll the
#    bottom pane so row 29 is at the top, and select the new D38 cell
#    (mirrors the source file's updated sheetView/selection).
# ---------------------------------------------------------------------
try {
    $win = $excel.ActiveWindow

    $win.FreezePanes = $false
    $ws.Range("A2").Select()
    $win.FreezePanes = $true

    $win.ScrollRow = 29
    $win.ScrollColumn = 1

    $ws.Range("D38").Select()
} catch {
    # View-state adjustments are best-effort; ignore if unsupported.
}
